$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Currículo – Alex Wilber", $true, $false, $false, $false, $false,
    $true, 1, $false, "Currículo: Alex Wilber", 2)

$d.Content.Find.Execute(
    "Spark Animation: Designer de Animação (Jan 2021 - Presente)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Spark Animation: Designer de animação (janeiro de 2021 – presente)", 2)

$d.Content.Find.Execute(
    "Pixel Studio: Designer de Animação (Jun 2018 - Dez 2020)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Pixel Studio: Designer de animação (jun. 2018 – dez. 2020)", 2)

$d.Content.Find.Execute(
    "Animação Flash: Designer de Animação Júnior (Set 2016 - Maio 2018)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Flash Animation: Designer de Animação Júnior (set 2016 – maio 2018)", 2)

$d.Content.Find.Execute(
    "Mestrado em Animação, Previsão de Formatura: Dez 2025", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mestrado em Animação, Conclusão esperada: Dec. 2025", 2)

$d.Content.Find.Execute(
    "A Arte da Animação 3D: Um Guia para Iniciantes.", $true, $false, $false, $false, $false,
    $true, 1, $false, "The Art of 3D Animation: A Guide for Beginners.", 2)

$d.Content.Find.Execute(
    "Nova Iorque: Spark Press.", $true, $false, $false, $false, $false,
    $true, 1, $false, "New York: Spark Press.", 2)
